$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# EPBDS-4097 Deployment isolation feature, fix bug in class loader isolation
#
# project2/Module2_1.xlsx gains a second "printJavaBeanSecond()" rule that
# mirrors the existing "printJavaBean()" rule, the "beans" data table is
# renamed to "beans2" (and all references to it updated), and the "import"
# value is changed from org.openl.example to org.openl.example2.
# ---------------------------------------------------------------------------

# Insert two new rows (18:19) right after the existing "printDatatype()"
# rule block (rows 15:16) to host the new "printJavaBeanSecond()" rule.
# This shifts the old rows 19:21 ("Environment" / "import" tables) down to
# 21:23, exactly as it happens when a user inserts rows in Excel.
$ws.Rows("18:19").Insert()

# Re-create the merged, bordered, centered look of the other rule-name rows
# (e.g. B16:D16) for the two freshly inserted rows.
$ws.Range("B18:D18").Merge()
$ws.Range("B19:D19").Merge()
$ws.Range("B18:D18").Borders.LineStyle = 1
$ws.Range("B19:D19").Borders.LineStyle = 1
$ws.Range("B18:D18").HorizontalAlignment = -4108
$ws.Range("B19:D19").HorizontalAlignment = -4108

# New rule: "Method String printJavaBeanSecond()" that prints beans2[0]
$ws.Range("B18").Value = "Method String printJavaBeanSecond()"
$ws.Range("B19").Value = "return ""project2""+print(beans2[0]);"

# The "beans" data table parameter is renamed to "beans2"...
$ws.Range("F4").Value = "Data TestBean beans2"

# ...so the existing printJavaBean() rule now references beans2 as well.
$ws.Range("B12").Value = "return ""project2""+print(beans2[0]);"

# The shared "import" table now points at org.openl.example2.
$ws.Range("C23").Value = "org.openl.example2"

# Restore the user's last cell selection as captured in the saved workbook.
$ws.Range("F4:G4").Select()
